$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row layout (columns A through N)
$ws.Range("A1").Value = "AreaMeasurementId"
$ws.Range("B1").Value = "Standard"
$ws.Range("C1").Value = "Type"
$ws.Range("D1").Value = "Value"
$ws.Range("E1").Value = "Unit"
$ws.Range("F1").Value = "ValidFrom"
$ws.Range("G1").Value = "ValidUntil"
$ws.Range("H1").Value = "SiteId"
$ws.Range("I1").Value = "LandId"
$ws.Range("J1").Value = "BuildingId"
$ws.Range("K1").Value = "UnitId"
$ws.Range("L1").Value = "SpaceId"
$ws.Range("M1").Value = "FloorId"
$ws.Range("N1").Value = "Guid"

# UnitId header is bold (new cell style with applyFont)
$ws.Range("K1").Font.Bold = $true

# Column A autofit width (bestFit), matches observed column width 17.6640625
$ws.Columns("A").AutoFit()

# Update selection to A4
$ws.Range("A4").Select()
